# Slide 18 ("For next week"): the lab announcement paragraphs are updated to
# mention that Lab 2 also went out, and the surrounding explanation text is
# reworded to match (commit: "sept9 rec with lab 2").
#
# Each paragraph's text is first swapped to a disjoint placeholder string
# before being set to its final value. PowerPoint's TextRange.Text setter
# otherwise tries to reuse/split runs around any text shared between the old
# and new values (e.g. the common "Lab"/"T" prefix), which would leave the
# run fragmented. Going through an unrelated placeholder string first keeps
# each paragraph a single clean <a:r> run, matching a normal in-place edit.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(18)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

$tr.Paragraphs(1,1).Text = "###PLACEHOLDER###"
$tr.Paragraphs(1,1).Text = "Labs 1 and 2 are out…!"

$tr.Paragraphs(2,1).Text = "###PLACEHOLDER###"
$tr.Paragraphs(2,1).Text = "They’re relatively easy – one is showing me you have Java installed and working, and another is answering some questions from lecture (open note!)"

$tr.Paragraphs(3,1).Text = "###PLACEHOLDER###"
$tr.Paragraphs(3,1).Text = "Keep an eye on Canvas for further submission instructions; it’ll be posted some time tonight (email me if Friday comes and I forgot to post the assignment). "
